$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4006.4285
$ws.Range("I51").Value = 3498.5
$ws.Range("K51").Value = 3498.5
$ws.Range("M51").Value = -3014.5

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 19122.666
$ws.Range("I62").Value = 14694
$ws.Range("J62").Value = 27980
$ws.Range("K62").Value = 14694
$ws.Range("L62").Value = 27980
$ws.Range("M62").Value = -14070
$ws.Range("N62").Value = -29228

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 19122.666
$ws.Range("I65").Value = 14694
$ws.Range("J65").Value = 27980
$ws.Range("K65").Value = 73470
$ws.Range("L65").Value = 139900
$ws.Range("M65").Value = -70350
$ws.Range("N65").Value = -146140

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4185.9443
$ws.Range("J138").Value = 2205.9
$ws.Range("L138").Value = 6617.700000000001
$ws.Range("N138").Value = -16897.7

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23255812
$ws.Range("I2").Value = 23255812
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 23255812
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -23255699
$ws.Range("N2").ClearContents()

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3277.5286
$ws.Range("I32").Value = 2778.4182
$ws.Range("K32").Value = 2778.4182
$ws.Range("M32").Value = -2491.4182

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3071.4285
$ws.Range("I45").Value = 4250
$ws.Range("J45").Value = 2794.1177
$ws.Range("K45").Value = 4250
$ws.Range("L45").Value = 2794.1177
$ws.Range("M45").Value = -3873
$ws.Range("N45").Value = -3548.1177

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2237.5862
$ws.Range("I61").Value = 1221.591
$ws.Range("K61").Value = 1221.591
$ws.Range("M61").Value = -1009.591

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 839.4643
$ws.Range("I74").Value = 508.6
$ws.Range("K74").Value = 508.6
$ws.Range("M74").Value = 365.4

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 839.4643
$ws.Range("I77").Value = 508.6
$ws.Range("K77").Value = 2543
$ws.Range("M77").Value = 1825

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 803.1667
$ws.Range("I110").Value = 803.1667
$ws.Range("K110").Value = 803.1667
$ws.Range("M110").Value = 1241.8333

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 23255812
$ws.Range("I116").Value = 23255812
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 23255812
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -23253518
$ws.Range("N116").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1877.95
$ws.Range("I132").Value = 1506.0667
$ws.Range("J132").Value = 2993.6
$ws.Range("K132").Value = 4518.2001
$ws.Range("L132").Value = 8980.799999999999
$ws.Range("M132").Value = -1988.2001
$ws.Range("N132").Value = -14040.8

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 100261
$ws.Range("J133").Value = 100261
$ws.Range("L133").Value = 100261
$ws.Range("N133").Value = -105321

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2237.5862
$ws.Range("I136").Value = 1221.591
$ws.Range("K136").Value = 3664.773
$ws.Range("M136").Value = -1114.773

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23255812
$ws.Range("I3").Value = 23255812
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 23255812
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -23255698
$ws.Range("N3").ClearContents()

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 199
$ws.Range("I22").Value = 199
$ws.Range("K22").Value = 199
$ws.Range("M22").Value = -26

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2044.1111
$ws.Range("I94").Value = 1987.125
$ws.Range("K94").Value = 1987.125
$ws.Range("M94").Value = -1536.125

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2992.182
$ws.Range("I134").Value = 3054.4375
$ws.Range("K134").Value = 9163.3125
$ws.Range("M134").Value = -6628.3125

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5209015
$ws.Range("I22").Value = 572.1667
$ws.Range("J22").Value = 10417458
$ws.Range("K22").Value = 572.1667
$ws.Range("L22").Value = 10417458
$ws.Range("M22").Value = -222.1667
$ws.Range("N22").Value = -10418158

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1113545.8
$ws.Range("I99").Value = 2002179.6
$ws.Range("J99").Value = 2753.5
$ws.Range("K99").Value = 2002179.6
$ws.Range("L99").Value = 2753.5
$ws.Range("M99").Value = -2000681.6
$ws.Range("N99").Value = -5749.5

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1113545.8
$ws.Range("I126").Value = 2002179.6
$ws.Range("J126").Value = 2753.5
$ws.Range("K126").Value = 6006538.800000001
$ws.Range("L126").Value = 8260.5
$ws.Range("M126").Value = -6004068.800000001
$ws.Range("N126").Value = -13200.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2083.7273
$ws.Range("I132").Value = 1636.2941
$ws.Range("K132").Value = 4908.8823
$ws.Range("M132").Value = -2378.8823

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1795.069
$ws.Range("I134").Value = 1034.1428
$ws.Range("K134").Value = 3102.4284
$ws.Range("M134").Value = -567.4284000000002

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 84332.664
$ws.Range("J37").Value = 84332.664
$ws.Range("L37").Value = 252997.992
$ws.Range("N37").Value = -253221.992

# CUL row 123
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 125001650
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1664.8
$ws.Range("I134").Value = 1664.8
$ws.Range("K134").Value = 4994.4
$ws.Range("M134").Value = 75.60000000000036

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5666.3335
$ws.Range("I102").Value = 7499.5
$ws.Range("K102").Value = 7499.5
$ws.Range("M102").Value = -5877.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2751236.8
$ws.Range("I132").Value = 4811178.5
$ws.Range("J132").Value = 4647.5
$ws.Range("K132").Value = 14433535.5
$ws.Range("L132").Value = 13942.5
$ws.Range("M132").Value = -14431005.5
$ws.Range("N132").Value = -19002.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5076.846
$ws.Range("I7").Value = 2857.5715
$ws.Range("J7").Value = 7666
$ws.Range("K7").Value = 2857.5715
$ws.Range("L7").Value = 7666
$ws.Range("M7").Value = -2745.5715
$ws.Range("N7").Value = -7890

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7912.5
$ws.Range("I40").Value = 6825.5
$ws.Range("J40").Value = 8999.5
$ws.Range("K40").Value = 6825.5
$ws.Range("L40").Value = 8999.5
$ws.Range("M40").Value = -6689.5
$ws.Range("N40").Value = -9271.5

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8688.375
$ws.Range("I122").Value = 6584.6665
$ws.Range("J122").Value = 14999.5
$ws.Range("K122").Value = 19753.9995
$ws.Range("L122").Value = 44998.5
$ws.Range("M122").Value = -17303.9995
$ws.Range("N122").Value = -49898.5

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5076.846
$ws.Range("I126").Value = 2857.5715
$ws.Range("J126").Value = 7666
$ws.Range("K126").Value = 8572.7145
$ws.Range("L126").Value = 22998
$ws.Range("M126").Value = -6102.7145
$ws.Range("N126").Value = -27938

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4446.409
$ws.Range("I132").Value = 997.5
$ws.Range("J132").Value = 4791.3
$ws.Range("K132").Value = 2992.5
$ws.Range("L132").Value = 14373.9
$ws.Range("M132").Value = -462.5
$ws.Range("N132").Value = -19433.9

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2192.4285
$ws.Range("I132").Value = 1087.5
$ws.Range("K132").Value = 3262.5
$ws.Range("M132").Value = -732.5
